$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header "all epoch incorrect guesses" in W1, matching the
# formatting of the neighbouring header cell V1.
$ws.Range("V1").Copy()
$ws.Range("W1").PasteSpecial(-4122)
$ws.Range("W1").Value = "all epoch incorrect guesses"

# Update the run's recorded results: training now actually happened, so
# "No training occurred" placeholders are replaced with real values.
$ws.Range("U2").Value = 0
$ws.Range("V2").Value = 100
$ws.Range("W2").Value = "[]"
